$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1 headers (new columns E-I), written left to right so the
#     shared-string table picks them up in that order ---
$ws.Range("E1").Value = "检验时间"
$ws.Range("F1").Value = "检验类型"
$ws.Range("G1").Value = "细菌名称"
$ws.Range("H1").Value = "药1"
$ws.Range("I1").Value = "药2"

# --- Row 2 new data; string values written in the same order the
#     shared-string table picked them up in the target file (H2 first) ---
$ws.Range("H2").Value = "耐药"
$ws.Range("F2").Value = "3,4"
$ws.Range("G2").Value = "牛逼"

# Date value mirrors the existing C2 cell formatting (built-in numFmtId 14)
$ws.Range("E2").NumberFormat = "m/d/yy"
$ws.Range("E2").Value = 41996

# Column width for the new column E (raw OOXML width of 22)
$ws.Range("E1").EntireColumn.ColumnWidth = 21.285714285714285

# Update the active selection to match the edited file
$ws.Range("G5").Select()
